$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates (Authors, Other found locations, Misc. Data)
$ws.Range("E2").Value = "[Jennifer%Lighter%Jennifer.Lighter@nyumc.org%0, Michael%Phillips%NULL%1, Sarah%Hochman%NULL%1, Stephanie%Sterling%NULL%1, Diane%Johnson%NULL%1, Fritz%Francois%NULL%0, Anna%Stachel%NULL%1]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Oxford University Press"

# Row 3 updates (Title, Authors, ID, ID Format, Date Accepted, Misc. Data)
$ws.Range("C3").Value = "Unknown Title"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "1970-01-01"
$ws.Range("H3").Style = "Normal"
$ws.Range("J3").Value = ""
